{"js": "// Replace the closing remark of the lab report with a more specific\n// sentence, splitting off the product name \"MatLab\" into its own run\n// (it picks up a distinct font attribute, mirroring the author's edit).\n\nconst body = context.document.body;\n\nconst oldText = \"\u041d\u0438\u0447\u0435\u0433\u043e \u043d\u0435 \u043f\u043e\u043d\u044f\u0442\u043d\u043e, \u043d\u043e  \u043e\u0447\u0435\u043d\u044c \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u043e\";\nconst newLead = \"\u041e\u0437\u043d\u0430\u043a\u043e\u043c\u0438\u043b\u0441\u044f \u0441 \u0440\u0430\u0431\u043e\u0442\u043e\u0439 \u043f\u0440\u043e\u0441\u0442\u0435\u0439\u0448\u0438\u0445 \u0447\u0438\u0441\u0435\u043b \u0432 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0435 \";\nconst newTail = \"MatLab\";\n\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \" + oldText);\n}\n\n// Replace the whole old sentence with the new lead-in text; the run keeps\n// its original formatting (color / size) because we are editing in place.\nconst target = results.items[0];\ntarget.insertText(newLead, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-locate the freshly inserted lead-in text so we can anchor the new\n// \"MatLab\" run immediately after it (inside the same paragraph, before the\n// trailing period).\nconst leadResults = body.search(newLead, { matchCase: true, matchWholeWord: false });\nleadResults.load(\"items\");\nawait context.sync();\n\nconst lead = leadResults.items[0];\nconst tailRange = lead.insertText(newTail, Word.InsertLocation.end);\n\n// Give the new run its own font attribute so it is emitted as a separate\n// <w:r> (matching the original document's run split).\ntailRange.font.set({ name: \"Calibri\" });\n\nawait context.sync();\n", "ps1": "# Replace the closing remark of the lab report with a more specific\n# sentence, splitting off the product name \"MatLab\" into its own run\n# (it picks up a distinct font attribute, mirroring the author's edit).\n\n$d = $word.ActiveDocument\n\n$oldText = \"\u041d\u0438\u0447\u0435\u0433\u043e \u043d\u0435 \u043f\u043e\u043d\u044f\u0442\u043d\u043e, \u043d\u043e  \u043e\u0447\u0435\u043d\u044c \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u043e\"\n$newLead = \"\u041e\u0437\u043d\u0430\u043a\u043e\u043c\u0438\u043b\u0441\u044f \u0441 \u0440\u0430\u0431\u043e\u0442\u043e\u0439 \u043f\u0440\u043e\u0441\u0442\u0435\u0439\u0448\u0438\u0445 \u0447\u0438\u0441\u0435\u043b \u0432 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0435 \"\n$newTail = \"MatLab\"\n\n$find = $d.Content.Find\n$find.Text = $oldText\n$find.Execute()\n\nif ($find.Found) {\n    $target = $d.Content\n    $target.Find.Execute($oldText)\n\n    # Replace the old sentence in place so the run keeps its original\n    # formatting (color / size).\n    $target.Text = $newLead\n\n    # Insert the new run immediately after the lead-in text, before the\n    # trailing period that already follows in the paragraph.\n    $insertPoint = $target.End\n    $tailRange = $d.Range($insertPoint, $insertPoint)\n    $tailRange.Text = $newTail\n\n    # Give the new run its own font attribute so it is emitted as a\n    # separate <w:r> (matching the original document's run split).\n    $tailRange = $d.Range($insertPoint, $insertPoint + $newTail.Length)\n    $tailRange.Font.Name = \"Calibri\"\n}\n"}
